# DOMA-11339 - add "Decommissioning date" column and trim the unused
# template filler rows from the meter-import example workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "Decommissioning date" column (U) -----------------------
# Header cell, with the same text/data styling as the neighbouring
# "Automatic" column (T).
$ws.Range("U1").Value = "Decommissioning date"

$ws.Range("T1").Copy()
$ws.Range("U1").PasteSpecial(-4122)

$ws.Range("T2").Copy()
$ws.Range("U2").PasteSpecial(-4122)

$ws.Range("U1").ColumnWidth = $ws.Range("T1").ColumnWidth

# --- Drop the empty styled template rows (3-10) ---------------------------
$ws.Range("A3:T10").EntireRow.Delete()
